$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand new row at position 391, pushing the existing rows
# 391..518 down to 392..519 (dimension grows from A1:R518 to A1:R519).
$ws.Rows.Item(391).Insert()

# Populate the newly inserted row 391 with the new record.
$ws.Cells.Item(391, 1).Value = 9
$ws.Cells.Item(391, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(391, 3).Value = "Metropolitana"
$ws.Cells.Item(391, 4).Value = 44855
$ws.Cells.Item(391, 5).Value = 13
$ws.Cells.Item(391, 6).Value = 100112028
$ws.Cells.Item(391, 7).Value = "Sandia"
$ws.Cells.Item(391, 8).Value = "Sin especificar"
$ws.Cells.Item(391, 9).Value = "Primera"
$ws.Cells.Item(391, 10).Value = 320
$ws.Cells.Item(391, 11).Value = 1000
$ws.Cells.Item(391, 12).Value = 1000
$ws.Cells.Item(391, 13).Value = 1000
$ws.Cells.Item(391, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(391, 15).Value = "Perú"
$ws.Cells.Item(391, 16).Value = 1000
$ws.Cells.Item(391, 17).Value = 1
$ws.Cells.Item(391, 18).Value = "Hortaliza"
